$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D4").Value = -7.877

$ws.Range("C10").Value = -12.948

$ws.Range("C12").Value = -11.466
$ws.Range("D12").Value = -7.775

$ws.Range("D17").Value = -8.282

$ws.Range("C18").Value = -12.406

$ws.Range("D26").Value = -7.255000000000001

$ws.Range("D27").Value = -7.685

$ws.Range("D28").Value = -7.890000000000001

$ws.Range("C37").Value = -13.179
$ws.Range("D37").Value = -7.369999999999999

$ws.Range("C55").Value = -13.633

$ws.Range("D65").Value = -7.67

$ws.Range("C68").Value = -10.809

$ws.Range("D73").Value = -7.826000000000001

$ws.Range("C77").Value = -13.169

$ws.Range("C78").Value = -13.004

$ws.Range("D84").Value = -8.242000000000001

$ws.Range("D85").Value = -8.563999999999998

$ws.Range("D93").Value = -7.007000000000001

$ws.Range("D95").Value = -7.569

$ws.Range("D98").Value = -7.231

$ws.Range("D99").Value = -8.193000000000001

$ws.Range("D101").Value = -7.825

$wb.Save()
